# Apply "output generated at 456a3b4" update to 广州-漫展信息.xlsx
# Updates the "想去人数" (want-to-go count, column F) figures that changed
# between crawls, refreshes two Cover image URLs, and appends a newly
# scraped concert ("广州·浪漫古典·百年经典世界名曲音乐会") to the 演出 sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (Exhibitions)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value2  = 7678
$ws.Range("F3").Value2  = 99
$ws.Range("F5").Value2  = 7159
$ws.Range("F8").Value2  = 623
$ws.Range("F11").Value2 = 436
$ws.Range("F12").Value2 = 769
$ws.Range("F13").Value2 = 33
$ws.Range("F14").Value2 = 72
$ws.Range("I14").Value2 = "//i2.hdslb.com/bfs/openplatform/202403/I6vXbS291711363371717.jpeg"
$ws.Range("F15").Value2 = 291
$ws.Range("F17").Value2 = 258
$ws.Range("F18").Value2 = 134
$ws.Range("F20").Value2 = 141
$ws.Range("F21").Value2 = 1085
$ws.Range("F23").Value2 = 596
$ws.Range("F24").Value2 = 2188
$ws.Range("F25").Value2 = 717
$ws.Range("F26").Value2 = 45
$ws.Range("F27").Value2 = 48
$ws.Range("F29").Value2 = 603
$ws.Range("F30").Value2 = 46

# ---------------------------------------------------------------------
# Sheet "演出" (Performances) - update F4, then append new row 10
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value2 = 320

# Copy the formatting of the existing index cell (A9, which carries the
# bold/centered/bordered style used throughout column A) down onto the
# new row's index cell before writing its value.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value2 = 9

# B10 holds a plain "yyyy-mm-dd" label that must stay literal text (not
# get auto-coerced into a date serial number). Temporarily force a text
# format so the assignment isn't reinterpreted as a date, then restore
# the cell's formatting (matching its neighbours, which carry the sheet's
# default/general style) now that the literal text is safely stored.
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value2 = "2024-05-25"
$ws.Range("C9").Copy()
$ws.Range("B10").PasteSpecial(-4122)

$ws.Range("C10").Value2 = "广州·浪漫古典·百年经典世界名曲音乐会"
$ws.Range("D10").Value2 = "东风中路299号 广州中山纪念堂"
$ws.Range("E10").Value2 = "2024.05.25 20:00-05.25 21:30"
$ws.Range("F10").Value2 = 0
$ws.Range("G10").Value2 = 75
$ws.Range("H10").Value2 = "https://show.bilibili.com/platform/detail.html?id=83327"
$ws.Range("I10").Value2 = "//i2.hdslb.com/bfs/openplatform/202403/uRWx5ZEu1711079544682.jpeg"

# ---------------------------------------------------------------------
# Sheet "本地生活" (Local Life)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value2 = 441

# ---------------------------------------------------------------------
# Sheet "全部类型" (All Types) - mirrors the other sheets' rows
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value2  = 441
$ws.Range("F3").Value2  = 7678
$ws.Range("F4").Value2  = 99
$ws.Range("F7").Value2  = 7159
$ws.Range("F10").Value2 = 623
$ws.Range("F14").Value2 = 436
$ws.Range("F15").Value2 = 320
$ws.Range("F18").Value2 = 769
$ws.Range("F19").Value2 = 33
$ws.Range("F20").Value2 = 72
$ws.Range("I20").Value2 = "//i2.hdslb.com/bfs/openplatform/202403/I6vXbS291711363371717.jpeg"
$ws.Range("F21").Value2 = 291
$ws.Range("F26").Value2 = 258
$ws.Range("F27").Value2 = 134
$ws.Range("F29").Value2 = 141
$ws.Range("F30").Value2 = 1085
$ws.Range("F32").Value2 = 596
$ws.Range("F33").Value2 = 2188
$ws.Range("F34").Value2 = 717
$ws.Range("F35").Value2 = 45
$ws.Range("F36").Value2 = 48
$ws.Range("F39").Value2 = 603
$ws.Range("F40").Value2 = 46
